$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) from 2023-11-13 (45243) to 2023-11-14 (45244)
# for rows 2 through 6.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45244
}
